$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    68 = "Yes"
    69 = "No"
    70 = "No"
    71 = "No"
    72 = "No"
    73 = "No"
    74 = "No"
    75 = "Yes"
    76 = "No"
    77 = "No"
    78 = "No"
    79 = "No"
    80 = "No"
    81 = "No"
    82 = "No"
    83 = "No"
    84 = "No"
    85 = "No"
    86 = "No"
    87 = "No"
    88 = "No"
    89 = "No"
    90 = "No"
    91 = "No"
    92 = "No"
    93 = "No"
    94 = "No"
    95 = "No"
    96 = "Yes"
    97 = "No"
    98 = "Yes"
    99 = "No"
    100 = "No"
    101 = "No"
    102 = "No"
    103 = "No"
    104 = "No"
    105 = "No"
    106 = "No"
    107 = "No"
    108 = "No"
    109 = "No"
    110 = "No"
    111 = "No"
    112 = "No"
    113 = "No"
    114 = "No"
    115 = "No"
    116 = "No"
    117 = "No"
    118 = "No"
    119 = "No"
    120 = "No"
    121 = "Yes"
    122 = "No"
    123 = "Yes"
    124 = "Yes"
    125 = "No"
    126 = "No"
    127 = "No"
    128 = "Yes"
    129 = "No"
    130 = "No"
    131 = "No"
    132 = "No"
    133 = "No"
    134 = "No"
    135 = "No"
    136 = "Yes"
    137 = "No"
    138 = "No"
    139 = "No"
    140 = "No"
    141 = "No"
    142 = "No"
    143 = "No"
    144 = "No"
    145 = "No"
    146 = "No"
    147 = "No"
    148 = "No"
    149 = "No"
    150 = "No"
    151 = "No"
    152 = "No"
    153 = "No"
    154 = "No"
    155 = "No"
    156 = "No"
    157 = "No"
    158 = "No"
    159 = "No"
    160 = "No"
    161 = "No"
    162 = "No"
    163 = "No"
    164 = "No"
    165 = "No"
    166 = "No"
    167 = "No"
    168 = "Yes"
    169 = "No"
    171 = "No"
    172 = "No"
    173 = "No"
    174 = "No"
    175 = "No"
    176 = "No"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 10).Value = $values[$row]
}
